# Apply the "Consent codes -> 2 codes" + "Removed Facility Name/Number" edit
# described in the commit message to the Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Delete the three extra "Consent Decision Code" value rows (rows 24-26:
#    "Consent Denied", "Inmate Never Seen", "Consent Not Obtained"), and fold
#    their meaning into a single summary cell on the "Consent Decision Code"
#    row itself.
$ws.Range("B23").Value = "Codes: Consent Granted; Consent Denied"
$ws.Rows("24:26").Delete()

# 2) Remove the "Booking Facility Name" / "Booking Facility Number" rows
#    from the Booking section entirely.
$ws.Rows("17:18").Delete()

Write-Output "edit applied"
